# Generate Report for Handback
# - Mark rows as handed back (status + handback datetime)
# - Add "Latest Target File" / "Latest Handback File" columns (F/G) mirroring
#   the existing "Source File Name" / "Latest Handoff File" columns, since the
#   handback is in sync with en-US (same files were sent back unchanged).

$wb = $excel.ActiveWorkbook

$langSheets = @("zh-cn", "de-de")

foreach ($langName in $langSheets) {
    $ws = $wb.Worksheets.Item($langName)

    # ---- Status column (C): "Ready for handoff" -> "Handed back: in sync with en-US"
    $ws.Range("C2").Value = "Handed back: in sync with en-US"
    $ws.Range("C3").Value = "Handed back: in sync with en-US"

    # ---- New column F: "Latest Target File" (mirrors column A, source file name)
    $aDisplay2 = $ws.Range("A2").Text
    $aAddress2 = $ws.Hyperlinks.Item(1).Address
    $ws.Range("F2").Value = $aDisplay2
    $ws.Range("F2").Style = "HyperLink"
    $ws.Hyperlinks.Add($ws.Range("F2"), $aAddress2, [Type]::Missing, [Type]::Missing, $aDisplay2) | Out-Null

    $aDisplay3 = $ws.Range("A3").Text
    $ws.Range("F3").Value = $aDisplay3
    $ws.Range("F3").Style = "HyperLink"
    $ws.Hyperlinks.Add($ws.Range("F3"), $aAddress2, [Type]::Missing, [Type]::Missing, $aDisplay3) | Out-Null

    # ---- New column G: "Latest Handback File" (mirrors column D, latest handoff file)
    $dDisplay2 = $ws.Range("D2").Text
    $dAddress2 = $null
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address -eq $ws.Range("D2").Address) {
            $dAddress2 = $hl.Address
        }
    }
    $ws.Range("G2").Value = $dDisplay2
    $ws.Range("G2").Style = "HyperLink"
    $ws.Hyperlinks.Add($ws.Range("G2"), $dAddress2, [Type]::Missing, [Type]::Missing, $dDisplay2) | Out-Null

    $dDisplay3 = $ws.Range("D3").Text
    $dAddress3 = $null
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address -eq $ws.Range("D3").Address) {
            $dAddress3 = $hl.Address
        }
    }
    $ws.Range("G3").Value = $dDisplay3
    $ws.Range("G3").Style = "HyperLink"
    $ws.Hyperlinks.Add($ws.Range("G3"), $dAddress3, [Type]::Missing, [Type]::Missing, $dDisplay3) | Out-Null
}

# ---- Handback datetimes (column H): placeholder "0001-01-01 00:00:00" -> real timestamps.
# zh-cn finished syncing first, de-de a few seconds later.
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("H2").Value = "2016-03-17 14:10:49"
$wsZh.Range("H3").Value = "2016-03-17 14:10:49"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("H2").Value = "2016-03-17 14:10:58"
$wsDe.Range("H3").Value = "2016-03-17 14:10:58"
